# CIERRE 16 OCT 2021
# Advance the payroll workbook from "SEMANA 41" (04-10 OCT 2021) to
# "SEMANA 42" (11-17 OCT 2021): update the week-header text, this week's
# pay figures, and scroll the sheet down to the new block.

$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week header (shared string). Every other "SEMANA ..." cell on the sheet
# (H9, B27, H27, B43, H43, B60) is a formula that ultimately points back to
# B9, so updating B9 alone ripples through all of them.
$ws.Range("B9").Value = "SEMANA   42  DEL    11      Al    17   DE   OCTUBRE          2021"

# This week's pay-slip figures (K6/K24/K41 are SUM formulas and recompute
# automatically).
$ws.Range("K4").Value = 1300
$ws.Range("K21").Value = 2660
$ws.Range("K40").Value = 0

# Scroll the saved view down to the new active block (was topLeftCell A40).
$win = $excel.ActiveWindow
$win.ScrollRow = 46
$win.ScrollColumn = 1
